$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 12 (hunk 0)
$ws.Range("H12").Value = 643.2308
$ws.Range("I12").Value = 626.0833
$ws.Range("J12").Value = 849
$ws.Range("K12").Value = 626.0833
$ws.Range("L12").Value = 849
$ws.Range("M12").Value = -456.0833
$ws.Range("N12").Value = -1189
# row 17 (hunk 1)
$ws.Range("H17").Value = 1471.2028
$ws.Range("J17").Value = 1471.2028
$ws.Range("L17").Value = 4413.6084
$ws.Range("N17").Value = -4749.6084
# row 31 (hunk 2)
$ws.Range("H31").Value = 370
$ws.Range("I31").Value = 370
$ws.Range("K31").Value = 1110
$ws.Range("M31").Value = -880
# row 74 (hunk 3)
$ws.Range("H74").Value = 3495
$ws.Range("I74").Value = 3495
$ws.Range("K74").Value = 3495
$ws.Range("M74").Value = -2559
# row 77 (hunk 4)
$ws.Range("H77").Value = 3495
$ws.Range("I77").Value = 3495
$ws.Range("K77").Value = 17475
$ws.Range("M77").Value = -12795
# row 100 (hunk 5)
$ws.Range("H100").Value = 2597.75
$ws.Range("I100").Value = 2427.4
$ws.Range("K100").Value = 2427.4
$ws.Range("M100").Value = -1886.4
# row 137 (hunk 6)
$ws.Range("H137").Value = 1976.5385
$ws.Range("I137").Value = 1538.409
$ws.Range("K137").Value = 4615.227000000001
$ws.Range("M137").Value = -2065.227000000001
# row 138 (hunk 7)
$ws.Range("H138").Value = 1709.7894
$ws.Range("I138").Value = 678.2143
$ws.Range("K138").Value = 2034.6429
$ws.Range("M138").Value = 3105.3571

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 45 (hunk 8)
$ws.Range("H45").Value = 2626.7144
$ws.Range("I45").Value = 2566.1667
$ws.Range("K45").Value = 2566.1667
$ws.Range("M45").Value = -2189.1667
# row 74 (hunk 9)
$ws.Range("H74").Value = 2729.6667
$ws.Range("I74").Value = 2776
$ws.Range("J74").Value = 2498
$ws.Range("K74").Value = 2776
$ws.Range("L74").Value = 2498
$ws.Range("M74").Value = -1902
$ws.Range("N74").Value = -4246
# row 77 (hunk 10)
$ws.Range("H77").Value = 2729.6667
$ws.Range("I77").Value = 2776
$ws.Range("J77").Value = 2498
$ws.Range("K77").Value = 13880
$ws.Range("L77").Value = 12490
$ws.Range("M77").Value = -9512
$ws.Range("N77").Value = -21226
# row 97 (hunk 11)
$ws.Range("H97").Value = 774
$ws.Range("I97").Value = 207.8
$ws.Range("J97").Value = 5020.5
$ws.Range("K97").Value = 207.8
$ws.Range("L97").Value = 5020.5
$ws.Range("M97").Value = 288.2
$ws.Range("N97").Value = -6012.5
# row 110 (hunk 12)
$ws.Range("H110").Value = 5097.364
$ws.Range("I110").Value = 5097.364
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 5097.364
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = -3052.364
$ws.Range("N110").ClearContents()
# row 122 (hunk 13)
$ws.Range("H122").Value = 1417.1562
$ws.Range("I122").Value = 1472
$ws.Range("K122").Value = 4416
$ws.Range("M122").Value = -1966
# row 134 (hunk 14)
$ws.Range("H134").Value = 70000
$ws.Range("J134").Value = 70000
$ws.Range("L134").Value = 70000
$ws.Range("N134").Value = -80140

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 59 (hunk 15)
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
# row 86 (hunk 16)
$ws.Range("H86").Value = 5917
$ws.Range("I86").Value = 5917
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 5917
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -4794
$ws.Range("N86").ClearContents()
# row 89 (hunk 17)
$ws.Range("H89").Value = 5917
$ws.Range("I89").Value = 5917
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 29585
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -23969
$ws.Range("N89").ClearContents()
# row 105 (hunk 18)
$ws.Range("H105").Value = 4547.7104
$ws.Range("I105").Value = 4015.76
$ws.Range("K105").Value = 4015.76
$ws.Range("M105").Value = -2268.76

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 7 (hunk 19)
$ws.Range("H7").Value = 345.91666
$ws.Range("I7").Value = 389.55554
$ws.Range("K7").Value = 389.55554
$ws.Range("M7").Value = -276.55554
# row 16 (hunk 20)
$ws.Range("H16").Value = 2469.7273
$ws.Range("I16").Value = 2533.6
$ws.Range("J16").Value = 2416.5
$ws.Range("K16").Value = 2533.6
$ws.Range("L16").Value = 2416.5
$ws.Range("M16").Value = -2246.6
$ws.Range("N16").Value = -2990.5
# row 31 (hunk 21)
$ws.Range("H31").Value = 1888.1428
$ws.Range("I31").Value = 1869.8334
$ws.Range("K31").Value = 1869.8334
$ws.Range("M31").Value = -1574.8334
# row 34 (hunk 22)
$ws.Range("H34").Value = 1888.1428
$ws.Range("I34").Value = 1869.8334
$ws.Range("K34").Value = 1869.8334
$ws.Range("M34").Value = -1667.8334
# row 62 (hunk 23)
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
# row 65 (hunk 24)
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
# row 99 (hunk 25)
$ws.Range("H99").Value = 3333.3333
$ws.Range("I99").Value = 3000
$ws.Range("K99").Value = 3000
$ws.Range("M99").Value = -1502
# row 113 (hunk 26)
$ws.Range("H113").Value = 2469.7273
$ws.Range("I113").Value = 2533.6
$ws.Range("J113").Value = 2416.5
$ws.Range("K113").Value = 2533.6
$ws.Range("L113").Value = 2416.5
$ws.Range("M113").Value = -363.5999999999999
$ws.Range("N113").Value = -6756.5
# row 122 (hunk 27)
$ws.Range("H122").Value = 7835.3335
$ws.Range("J122").Value = 3004.3333
$ws.Range("L122").Value = 9012.999899999999
$ws.Range("N122").Value = -13912.9999
# row 126 (hunk 28)
$ws.Range("H126").Value = 3333.3333
$ws.Range("I126").Value = 3000
$ws.Range("K126").Value = 9000
$ws.Range("M126").Value = -6530

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 2 (hunk 29)
$ws.Range("H2").Value = 90.69231000000001
$ws.Range("I2").Value = 77.166664
$ws.Range("J2").Value = 102.28571
$ws.Range("K2").Value = 462.999984
$ws.Range("L2").Value = 613.71426
$ws.Range("M2").Value = -349.999984
$ws.Range("N2").Value = -839.71426
# row 4 (hunk 30)
$ws.Range("H4").Value = 1235169.8
$ws.Range("I4").Value = 44960.21
$ws.Range("K4").Value = 134880.63
$ws.Range("M4").Value = -134768.63
# row 34 (hunk 31)
$ws.Range("H34").Value = 4078.5833
$ws.Range("I34").Value = 583.6667
$ws.Range("J34").Value = 5243.5557
$ws.Range("K34").Value = 1751.0001
$ws.Range("L34").Value = 15730.6671
$ws.Range("M34").Value = -1667.0001
$ws.Range("N34").Value = -15898.6671
# row 46 (hunk 32)
$ws.Range("H46").Value = 632
$ws.Range("I46").Value = 632
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 1896
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -1805
$ws.Range("N46").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 113 (hunk 33)
$ws.Range("H113").Value = 2203.8572
$ws.Range("I113").Value = 1986.4
$ws.Range("J113").Value = 2747.5
$ws.Range("K113").Value = 1986.4
$ws.Range("L113").Value = 2747.5
$ws.Range("M113").Value = 183.5999999999999
$ws.Range("N113").Value = -7087.5
# row 122 (hunk 34)
$ws.Range("H122").Value = 2639.6843
$ws.Range("I122").Value = 2157.6667
$ws.Range("J122").Value = 4447.25
$ws.Range("K122").Value = 6473.000100000001
$ws.Range("L122").Value = 13341.75
$ws.Range("M122").Value = -4023.000100000001
$ws.Range("N122").Value = -18241.75
# row 126 (hunk 35)
$ws.Range("H126").Value = 2971
$ws.Range("I126").Value = 2982.8333
$ws.Range("J126").Value = 2900
$ws.Range("K126").Value = 8948.499899999999
$ws.Range("L126").Value = 8700
$ws.Range("M126").Value = -6478.499899999999
$ws.Range("N126").Value = -13640
# row 132 (hunk 36)
$ws.Range("H132").Value = 2975
$ws.Range("I132").Value = 2975
$ws.Range("K132").Value = 8925
$ws.Range("M132").Value = -6395

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 7 (hunk 37)
$ws.Range("H7").Value = 7155.9165
$ws.Range("I7").Value = 7197.4
$ws.Range("K7").Value = 7197.4
$ws.Range("M7").Value = -7085.4
# row 122 (hunk 38)
$ws.Range("H122").Value = 4064.9473
$ws.Range("I122").Value = 3653.0833
$ws.Range("J122").Value = 4771
$ws.Range("K122").Value = 10959.2499
$ws.Range("L122").Value = 14313
$ws.Range("M122").Value = -8509.249899999999
$ws.Range("N122").Value = -19213
# row 126 (hunk 39)
$ws.Range("H126").Value = 7155.9165
$ws.Range("I126").Value = 7197.4
$ws.Range("K126").Value = 21592.2
$ws.Range("M126").Value = -19122.2

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 46 (hunk 40)
$ws.Range("H46").Value = 28429
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 28429
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 28429
$ws.Range("N46").Value = -28891
$ws.Range("M46").ClearContents()
# row 61 (hunk 41)
$ws.Range("H61").Value = 28579.8
$ws.Range("I61").Value = 7499.5
$ws.Range("K61").Value = 7499.5
$ws.Range("M61").Value = -7207.5
# row 107 (hunk 42)
$ws.Range("H107").Value = 347.75
$ws.Range("I107").Value = 347.75
$ws.Range("K107").Value = 1043.25
$ws.Range("M107").Value = 876.75
# row 126 (hunk 43)
$ws.Range("H126").Value = 3323.3333
$ws.Range("I126").Value = 3368
$ws.Range("J126").Value = 3234
$ws.Range("K126").Value = 10104
$ws.Range("L126").Value = 9702
$ws.Range("M126").Value = -7634
$ws.Range("N126").Value = -14642
# row 134 (hunk 44)
$ws.Range("H134").Value = 28429
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 28429
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 85287
$ws.Range("N134").Value = -90357
$ws.Range("M134").ClearContents()
# row 138 (hunk 45)
$ws.Range("H138").Value = 99997.5
$ws.Range("J138").Value = 99997.5
$ws.Range("L138").Value = 99997.5
$ws.Range("N138").Value = -110277.5
